$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Coin (B), Link (C), Price (D) and Volume(1h) (E) columns.
# For D-column values that look numeric, force text format first so Excel
# stores them as text (preserving exact formatting like trailing zeros)
# instead of silently converting them to numbers.

$ws.Range("D2").Value = '79.289.63'
$ws.Range("E2").Value = '  +3.96%  '
$ws.Range("D3").Value = '3.129.13'
$ws.Range("E3").Value = '  +1.73%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.997'
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '204.44'
$ws.Range("E5").Value = '  +2.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '618.89'
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.261'
$ws.Range("E7").Value = '  +23.79%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.998'
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.579'
$ws.Range("E9").Value = '  +5.02%  '
$ws.Range("D10").Value = '3.116.57'
$ws.Range("E10").Value = '  +1.34%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.574'
$ws.Range("E11").Value = '  +28.71%  '
$ws.Range("B12").Value = 'ShibaInu'
$ws.Range("C12").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000243'
$ws.Range("E12").Value = '  +24.46%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.163'
$ws.Range("E13").Value = '  +1.39%  '
$ws.Range("D14").Value = '3.676.51'
$ws.Range("E14").Value = '  +1.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.17'
$ws.Range("E15").Value = '  -1.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '31.01'
$ws.Range("E16").Value = '  +5.98%  '
$ws.Range("D17").Value = '79.203.97'
$ws.Range("E17").Value = '  +3.95%  '
$ws.Range("D18").Value = '3.117.50'
$ws.Range("E18").Value = '  +1.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.00'
$ws.Range("E19").Value = '  +3.38%  '
$ws.Range("B20").Value = 'SuiNetwork'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.90'
$ws.Range("E20").Value = '  +12.86%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '429.95'
$ws.Range("E21").Value = '  +12.52%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.95'
$ws.Range("E22").Value = '  -0.38%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.12'
$ws.Range("E23").Value = '  +14.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.76'
$ws.Range("E24").Value = '  +4.97%  '
$ws.Range("D25").Value = '3.282.31'
$ws.Range("E25").Value = '  +1.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '74.90'
$ws.Range("E26").Value = '  +3.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.60'
$ws.Range("E27").Value = '  +1.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.65'
$ws.Range("E28").Value = '  +6.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0000119'
$ws.Range("E30").Value = '  +9.88%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.996'
$ws.Range("E31").Value = '  -0.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.80'
$ws.Range("E32").Value = '  +5.99%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '546.15'
$ws.Range("E33").Value = '  +8.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.45'
$ws.Range("E34").Value = '  +1.26%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.147'
$ws.Range("E35").Value = '  +17.23%  '
$ws.Range("B36").Value = 'PancakeSwap'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.96'
$ws.Range("E36").Value = '  +1.87%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '22.56'
$ws.Range("E37").Value = '  +8.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.121'
$ws.Range("E38").Value = '  +18.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.998'
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.397'
$ws.Range("E40").Value = '  +5.02%  '
$ws.Range("B41").Value = 'WhiteBITCoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '20.69'
$ws.Range("E41").Value = '  +3.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '162.42'
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.47'
$ws.Range("E44").Value = '  +5.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '185.73'
$ws.Range("E45").Value = '  -4.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.77'
$ws.Range("E46").Value = '  +6.86%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.65'
$ws.Range("E47").Value = '  +8.10%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.764'
$ws.Range("E48").Value = '  -5.43%  '
$ws.Range("B49").Value = 'ImmutableX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.28'
$ws.Range("E49").Value = '  +1.65%  '
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '42.41'
$ws.Range("E50").Value = '  +4.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.16'
$ws.Range("E51").Value = '  +6.33%  '
